# "Generate Report for Archive"
#
# The localization status report is being refreshed: every "Ready for
# handoff" status cell becomes "In Translation", and the (now narrower)
# Status columns are resized to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    Overview!E2:F4 (per-language status columns) and the Status column
#    (C2:C4) on each per-language detail sheet all carry this value.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C4").Value = "In Translation"

# ---------------------------------------------------------------------
# 2) Narrow the Status columns (they no longer need to fit the longer
#    "Ready for handoff" text) from ~17.22 chars down to ~13.41 chars.
# ---------------------------------------------------------------------
$newStatusWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newStatusWidth
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth

$zhcn.Columns.Item(3).ColumnWidth = $newStatusWidth

$dede.Columns.Item(3).ColumnWidth = $newStatusWidth
